$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The '%syst_c' column (originally column O) is no longer used; remove the
# entire column so everything to its right shifts one place to the left.
$ws.Columns.Item(15).Delete()

# Rename the stat_plus/stat_minus/stat_u headers (now columns L, M, N) to
# tot_plus/tot_minus/tot_u.
$ws.Range("L1").Value = "tot_plus"
$ws.Range("M1").Value = "tot_minus"
$ws.Range("N1").Value = "tot_u"

# Match the saved selection state of the edited workbook.
$ws.Range("N2").Select() | Out-Null
